$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new row at row 32, shifting existing row 32+ down by one.
$ws.Rows.Item(32).Insert()

# Populate the new row 32 with the reward-points line item.
$ws.Cells.Item(32, 1).Value = "`$15 Off (300 points)"
$ws.Cells.Item(32, 37).Value = 300            # AK32
$ws.Cells.Item(32, 38).Value = "points"        # AL32
